$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 59560.53
$ws.Range("I28").Value = 67402
$ws.Range("K28").Value = 67402
$ws.Range("M28").Value = -66917
# Row 43
$ws.Range("H43").Value = 55557140
$ws.Range("I43").Value = 83334340
$ws.Range("J43").Value = 2750
$ws.Range("K43").Value = 83334340
$ws.Range("L43").Value = 2750
$ws.Range("M43").Value = -83334271
$ws.Range("N43").Value = -2888
# Row 62
$ws.Range("H62").Value = 3748
$ws.Range("I62").Value = 3748
$ws.Range("K62").Value = 3748
$ws.Range("M62").Value = -3124
# Row 65
$ws.Range("H65").Value = 3748
$ws.Range("I65").Value = 3748
$ws.Range("K65").Value = 18740
$ws.Range("M65").Value = -15620
# Row 98
$ws.Range("H98").Value = 2989.25
$ws.Range("I98").Value = 2340.5881
$ws.Range("K98").Value = 2340.5881
$ws.Range("M98").Value = -842.5880999999999
# Row 100
$ws.Range("H100").Value = 1358
$ws.Range("I100").Value = 1398
$ws.Range("J100").Value = 1198
$ws.Range("K100").Value = 1398
$ws.Range("L100").Value = 1198
$ws.Range("M100").Value = -857
$ws.Range("N100").Value = -2280
# Row 106
$ws.Range("H106").Value = 58826016
$ws.Range("I106").Value = 62502520
$ws.Range("K106").Value = 62502520
$ws.Range("M106").Value = -62501889
# Row 112
$ws.Range("H112").Value = 2985.9556
$ws.Range("J112").Value = 3084.9023
$ws.Range("L112").Value = 9254.706900000001
$ws.Range("N112").Value = -11470.7069
# Row 122
$ws.Range("H122").Value = 2989.25
$ws.Range("I122").Value = 2340.5881
$ws.Range("K122").Value = 7021.7643
$ws.Range("M122").Value = -4571.7643
# Row 132
$ws.Range("H132").Value = 1924.0358
$ws.Range("I132").Value = 1685.1818
$ws.Range("J132").Value = 2799.8333
$ws.Range("K132").Value = 5055.5454
$ws.Range("L132").Value = 8399.499899999999
$ws.Range("M132").Value = -2525.5454
$ws.Range("N132").Value = -13459.4999
# Row 137
$ws.Range("H137").Value = 1701694.1
$ws.Range("I137").Value = 5214.8535
$ws.Range("K137").Value = 15644.5605
$ws.Range("M137").Value = -13094.5605
# Row 141
$ws.Range("H141").Value = 4705.4375
$ws.Range("I141").Value = 4939.1333
$ws.Range("K141").Value = 14817.3999
$ws.Range("M141").Value = -9637.3999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 214128.72
$ws.Range("I32").Value = 238598.1
$ws.Range("K32").Value = 238598.1
$ws.Range("M32").Value = -238311.1
# Row 52
$ws.Range("H52").Value = 29998.334
$ws.Range("J52").Value = 29998.334
$ws.Range("L52").Value = 29998.334
$ws.Range("N52").Value = -30634.334
# Row 61
$ws.Range("H61").Value = 941161.9
$ws.Range("I61").Value = 25854.396
$ws.Range("K61").Value = 25854.396
$ws.Range("M61").Value = -25642.396
# Row 88
$ws.Range("H88").Value = 2869.7144
$ws.Range("J88").Value = 2865.1667
$ws.Range("L88").Value = 2865.1667
$ws.Range("N88").Value = -3677.1667
# Row 91
$ws.Range("H91").Value = 2869.7144
$ws.Range("J91").Value = 2865.1667
$ws.Range("L91").Value = 2865.1667
$ws.Range("N91").Value = -5673.1667
# Row 102
$ws.Range("H102").Value = 2088
$ws.Range("I102").Value = 1901.75
$ws.Range("J102").Value = 3205.5
$ws.Range("K102").Value = 1901.75
$ws.Range("L102").Value = 3205.5
$ws.Range("M102").Value = -279.75
$ws.Range("N102").Value = -6449.5
# Row 136
$ws.Range("H136").Value = 941161.9
$ws.Range("I136").Value = 25854.396
$ws.Range("K136").Value = 77563.18799999999
$ws.Range("M136").Value = -75013.18799999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3954.762
$ws.Range("I86").Value = 2296.4614
$ws.Range("K86").Value = 2296.4614
$ws.Range("M86").Value = -1173.4614
# Row 89
$ws.Range("H89").Value = 3954.762
$ws.Range("I89").Value = 2296.4614
$ws.Range("K89").Value = 11482.307
$ws.Range("M89").Value = -5866.307000000001
# Row 105
$ws.Range("H105").Value = 7025.1377
$ws.Range("I105").Value = 3102.3333
$ws.Range("J105").Value = 13444.272
$ws.Range("K105").Value = 3102.3333
$ws.Range("L105").Value = 13444.272
$ws.Range("M105").Value = -1355.3333
$ws.Range("N105").Value = -16938.272

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2681.2778
$ws.Range("I31").Value = 2904.3809
$ws.Range("J31").Value = 2368.9333
$ws.Range("K31").Value = 2904.3809
$ws.Range("L31").Value = 2368.9333
$ws.Range("M31").Value = -2609.3809
$ws.Range("N31").Value = -2958.9333
# Row 34
$ws.Range("H34").Value = 2681.2778
$ws.Range("I34").Value = 2904.3809
$ws.Range("J34").Value = 2368.9333
$ws.Range("K34").Value = 2904.3809
$ws.Range("L34").Value = 2368.9333
$ws.Range("M34").Value = -2702.3809
$ws.Range("N34").Value = -2772.9333
# Row 58
$ws.Range("H58").Value = 1642.6428
$ws.Range("J58").Value = 1603.6842
$ws.Range("L58").Value = 1603.6842
$ws.Range("N58").Value = -2009.6842
# Row 105
$ws.Range("H105").Value = 3116.6667
$ws.Range("I105").Value = 2633.3333
$ws.Range("K105").Value = 2633.3333
$ws.Range("M105").Value = -886.3332999999998
# Row 134
$ws.Range("H134").Value = 2446.1516
$ws.Range("I134").Value = 2252.7273
$ws.Range("K134").Value = 6758.1819
$ws.Range("M134").Value = -4223.1819
# Row 136
$ws.Range("H136").Value = 1642.6428
$ws.Range("J136").Value = 1603.6842
$ws.Range("L136").Value = 4811.0526
$ws.Range("N136").Value = -9911.052599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 276.33334
$ws.Range("I7").Value = 276.33334
$ws.Range("K7").Value = 829.0000200000001
$ws.Range("M7").Value = -717.0000200000001
# Row 11
$ws.Range("H11").Value = 454.16666
$ws.Range("I11").Value = 454.16666
$ws.Range("K11").Value = 1362.49998
$ws.Range("M11").Value = -1222.49998
# Row 80
$ws.Range("H80").Value = 807.25
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 807.25
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 92
$ws.Range("H92").Value = 531.6667
$ws.Range("I92").Value = 600
$ws.Range("K92").Value = 1800
$ws.Range("M92").Value = -552
# Row 107
$ws.Range("H107").Value = 945.4286
$ws.Range("J107").Value = 916.12
$ws.Range("L107").Value = 2748.36
$ws.Range("N107").Value = -6588.360000000001
# Row 136
$ws.Range("H136").Value = 7974.6333
$ws.Range("I136").Value = 4013.5557
$ws.Range("K136").Value = 12040.6671
$ws.Range("M136").Value = -6940.667099999999
# Row 138
$ws.Range("H138").Value = 4457.476
$ws.Range("I138").Value = 4438
$ws.Range("K138").Value = 13314
$ws.Range("M138").Value = -8174

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 714.6316
$ws.Range("I97").Value = 532.5833
$ws.Range("K97").Value = 532.5833
$ws.Range("M97").Value = -36.58330000000001
# Row 132
$ws.Range("H132").Value = 1283469.2
$ws.Range("I132").Value = 1320
$ws.Range("K132").Value = 3960
$ws.Range("M132").Value = -1430

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 3023824.5
$ws.Range("J43").Value = 3393125
$ws.Range("L43").Value = 3393125
$ws.Range("N43").Value = -3393511
# Row 55
$ws.Range("H55").Value = 559.78125
$ws.Range("I55").Value = 293.73334
$ws.Range("K55").Value = 293.73334
$ws.Range("M55").Value = -120.73334
# Row 82
$ws.Range("H82").Value = 1734
$ws.Range("I82").Value = 1681
$ws.Range("K82").Value = 1681
$ws.Range("M82").Value = -1320
# Row 85
$ws.Range("H85").Value = 1734
$ws.Range("I85").Value = 1681
$ws.Range("K85").Value = 1681
$ws.Range("M85").Value = -433
# Row 133
$ws.Range("H133").Value = 77777
$ws.Range("J133").Value = 77777
$ws.Range("L133").Value = 77777
$ws.Range("N133").Value = -82837
# Row 136
$ws.Range("H136").Value = 36151.4
$ws.Range("I136").Value = 86009
$ws.Range("J136").Value = 2913
$ws.Range("K136").Value = 258027
$ws.Range("L136").Value = 8739
$ws.Range("M136").Value = -255477
$ws.Range("N136").Value = -13839

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# Row 81
$ws.Range("H81").Value = 93093
$ws.Range("I81").Value = 2642.4443
$ws.Range("K81").Value = 5284.8886
$ws.Range("M81").Value = -4223.8886
# Row 84
$ws.Range("H84").Value = 93093
$ws.Range("I84").Value = 2642.4443
$ws.Range("K84").Value = 26424.443
$ws.Range("M84").Value = -21120.443
# Row 100
$ws.Range("H100").Value = 783.35297
$ws.Range("I100").Value = 582.625
$ws.Range("J100").Value = 3995
$ws.Range("K100").Value = 1165.25
$ws.Range("L100").Value = 7990
$ws.Range("M100").Value = -624.25
$ws.Range("N100").Value = -9072
# Row 132
$ws.Range("H132").Value = 66669064
$ws.Range("I132").Value = 142858980
$ws.Range("K132").Value = 428576940
$ws.Range("M132").Value = -428574410
# Row 136
$ws.Range("H136").Value = 1374.75
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1374.75
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 4124.25
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -9224.25
